$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 21: date, question text, and hyperlinked URL for "Coin Change"
$ws.Cells.Item(21, 1).Value = 46081
$ws.Cells.Item(21, 1).NumberFormat = "m/d/yyyy"
$ws.Cells.Item(21, 2).Value = "Coin Change"
$ws.Cells.Item(21, 3).Value = "https://leetcode.com/problems/coin-change/description/"
$ws.Hyperlinks.Add($ws.Cells.Item(21, 3), "https://leetcode.com/problems/coin-change/description/")
